$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 previously held A4="bag" + an already-empty B4. Drop the "bag" text
# (the cell itself goes away, not just its value) but keep B4 as-is.
$ws.Range("A4").Clear()

# Rows 5 ("saree") and 6 ("kurtis") are removed outright.
$ws.Range("A5:A6").EntireRow.Delete()

# New second-column values for rows 1-3.
$ws.Range("B1").Value = "print"
$ws.Range("B2").Value = "testcase2"
$ws.Range("B3").Value = "testcase3"

# Match the resulting selection/active cell recorded in the sheet view.
$ws.Range("A4:A6").Select()
